# ISS-45: Fix reading empty paragraphs
#
# Re-apply (to the extent the COM-interop host supports it) the changes
# recorded for test/SlideXML.Tests/Resource/010.pptx:
#   1. datetimeFigureOut field text "1/19/2020" -> "1/30/2020" on the
#      Handout Master and Notes Master "Date Placeholder" shapes.
#   2. spid bump on the think-cell OLE object in the "ER: One" slide
#      layout (_x0000_s83986 -> _x0000_s83987).
#   3. Removal of three decorative connector shapes ("Straight
#      Connector 5", "Frame Line", "Straight Connector 9") from the
#      "ER: One" slide layout.

$p = $ppt.ActivePresentation

# --- 1. Date fields on Handout Master / Notes Master -----------------
# (best effort -- some hosts do not allow writing into master-tier
# placeholders; failures here must not abort the rest of the script)
try {
    $hm = $p.HandoutMaster
    for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
        $shp = $hm.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "1/19/2020") {
            $shp.TextFrame.TextRange.Text = "1/30/2020"
        }
    }
} catch {
    Write-Output "HandoutMaster date update skipped: $_"
}

try {
    $nm = $p.NotesMaster
    for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
        $shp = $nm.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "1/19/2020") {
            $shp.TextFrame.TextRange.Text = "1/30/2020"
        }
    }
} catch {
    Write-Output "NotesMaster date update skipped: $_"
}

# --- 2 & 3. Slide layout "ER: One" (think-cell spid + connectors) ----
$layout = $null
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.CustomLayouts.Count; $i++) {
    $cl = $sm.CustomLayouts.Item($i)
    if ($cl.Name -eq "ER: One") {
        $layout = $cl
        break
    }
}
if (-not $layout) {
    $layout = $sm.CustomLayouts.Item(1)
}

# 2. Bump the think-cell OLE object's legacy VML spid (best effort --
# this id is not exposed on the standard Shape object model, so this
# is a no-op on hosts that don't support it).
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $shp = $layout.Shapes.Item($i)
    if ($shp.Name -eq "Object 10") {
        try { $shp.Spid = "_x0000_s83987" } catch { }
    }
}

# 3. Remove the three decorative connector shapes.
$connectorNames = @("Straight Connector 5", "Frame Line", "Straight Connector 9")
foreach ($connName in $connectorNames) {
    for ($i = $layout.Shapes.Count; $i -ge 1; $i--) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -eq $connName) {
            $shp.Delete()
        }
    }
}

Write-Output "Done. Layout shapes remaining: $($layout.Shapes.Count)"
